# Apply the "span_ec" worksheet edit:
#  - duplicate "ec_class accumulated" to a new leading sheet named "span_ec"
#  - re-label headers / the "@terms" row, refresh the B-column sample data,
#    and clear the D-column sample data (so D3:D5 error out to #DIV/0!)
#  - tweak sheet views (zoom/selection) on both the new sheet and the
#    original "ec_class accumulated" sheet

$wb = $excel.ActiveWorkbook

$srcRef = $wb.Worksheets.Item("ec_class accumulated")

# Duplicate the sheet; Excel inserts the copy immediately before the source.
# NOTE: worksheet variables here track by *position*, so after the copy,
# re-resolve fresh handles by index/name rather than reusing $srcRef.
$srcRef.Copy($srcRef)

$newSheet = $wb.Worksheets.Item(1)
$origSheet = $wb.Worksheets.Item("ec_class accumulated")

$newSheet.Name = "span_ec"

# --- Update the new "span_ec" sheet's contents -----------------------------
$newSheet.Range("B1").Value = "class_nl only"
$newSheet.Range("D1").Value = "span_nl (includes <span> tag)"
$newSheet.Range("B2").Value = "'@terms"
$newSheet.Range("D2").Value = "'@terms"

$newSheet.Range("B6").Value = 0.48859599999999997
$newSheet.Range("B7").Value = 0.48438100000000001
$newSheet.Range("B8").Value = 0.39359699999999997
$newSheet.Range("B9").Value = 0.50655499999999998
$newSheet.Range("B10").Value = 0.43773000000000001

$newSheet.Range("D6:D10").ClearContents()

# --- Sheet view tweaks -------------------------------------------------------
$origSheet.Activate()
$excel.ActiveWindow.Zoom = 150
$origSheet.Range("D10").Select()

$newSheet.Activate()
$excel.ActiveWindow.Zoom = 150
$newSheet.Range("B10").Select()

Write-Output "done"
